$d = $word.ActiveDocument

# Change 1: "AO JUÍZO DA {{ juizo_competente }}" -> "AO JUÍZO DA {{ juizo_competente.upper() }}"
$d.Content.Find.Execute("AO JUÍZO DA {{ juizo_competente }}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "AO JUÍZO DA {{ juizo_competente.upper() }}", 2)

# Change 2: "no que tange à ação  {{ " -> "no que tange à {{ "
$d.Content.Find.Execute("no que tange à ação  {{ ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "no que tange à {{ ", 2)
